# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet (3rd sheet) gains a new "category" column,
# inserted right after "property_category" (pushing the existing "date",
# "legislator_name" and "legislator_id" columns one position to the
# right), plus two brand-new trailing columns: "source_file" and "index".
#
# Final header/data layout for the stock sheet:
#   B name                 B2 中日國際企業股份有限公司
#   C owner                C2 林岱樺
#   D quantity              D2 500
#   E face_value           E2 10
#   F currency             F2 新臺幣
#   G total                G2 5000
#   H property_category    H2 stock
#   I category   (new)     I2 normal      (new)
#   J date                 J2 2012-02-29
#   K legislator_name      K2 林岱樺
#   L legislator_id        L2 904
#   M source_file (new)    M2 tmp3bff1    (new)
#   N index      (new)     N2 56          (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column at I; this shifts the existing I/J/K columns
# (date, legislator_name, legislator_id) right to J/K/L, carrying their
# header/data cell styles (s="1" / s="2") along automatically.
$ws.Columns("I:I").Insert()

# Grow two more trailing columns (M, N) by duplicating the now-rightmost
# data column (L) and inserting the copies right after it - this brings
# along the same header/data styles so the new cells match their
# neighbours instead of being left unstyled.
$ws.Columns("L:L").Copy()
$ws.Columns("M:M").Insert()

$ws.Columns("L:L").Copy()
$ws.Columns("N:N").Insert()

# New header row values.
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# New data row values.
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp3bff1"
$ws.Range("N2").Value = 56
